$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.31861424446106
$ws.Range("B1").Value = 2.437829494476318
$ws.Range("C1").Value = 5.855088233947754
$ws.Range("D1").Value = 1.724735021591187
$ws.Range("E1").Value = 1.291378617286682
